$d = $word.ActiveDocument

# Locate the last "Changes from Beta" bullet item - the paragraph that ends
# with "...add headers as suggested." - and append a new bullet list item
# after it (and before the trailing blank paragraphs), describing the new
# health system feature.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tidy up the variables*") {
        $target = $p
    }
}

$newPara = $target.Range.InsertParagraphAfter()

# Re-fetch the freshly inserted paragraph (it inherited the ListParagraph /
# numbered-list formatting from the paragraph it was split from) and give it
# its text.
$newRange = $target.Next().Range
$newRange.Text = "One of the peer reviewers suggested health for colliding with small asteroids – I have implemented a health system where the player starts with a maximum of 4 UFOs (varies depending on level), crashing into asteroids reduces this, and once it is 0 it’s game over.  Small asteroid removes 1, medium removes 2, large removes 3."
